$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 523.5861948266684
$ws.Range("D2").Value = 122.0436772156639
$ws.Range("G2").Value = 477
$ws.Range("H2").Value = 567
$ws.Range("C3").Value = 36.84388973146382
$ws.Range("D3").Value = 6.489012600754923
$ws.Range("F3").Value = 32.05
$ws.Range("G3").Value = 37.09
$ws.Range("H3").Value = 41.22
$ws.Range("C4").Value = 2.083656679063995
$ws.Range("D4").Value = 2.57745453201037
$ws.Range("F4").Value = 0.71
$ws.Range("G4").Value = 1.38
$ws.Range("H4").Value = 2.56
$ws.Range("C5").Value = 322.7980935511237
$ws.Range("D5").Value = 9.021253159600393
$ws.Range("F5").Value = 317.87
$ws.Range("G5").Value = 323.26
$ws.Range("H5").Value = 329.3
$ws.Range("C6").Value = 23.26879912621351
$ws.Range("D6").Value = 3.665981481259716
$ws.Range("F6").Value = 20.68
$ws.Range("G6").Value = 22.71
$ws.Range("H6").Value = 25.66
$ws.Range("C7").Value = -75.98254104315876
$ws.Range("D7").Value = 22.50803120461898
$ws.Range("C8").Value = 7.664853278133452
$ws.Range("D8").Value = 6.848737824671677
$ws.Range("C9").Value = 9.231929931115822
$ws.Range("D9").Value = 1.653348633136656
$ws.Range("C10").Value = 867.8271661262509
$ws.Range("D10").Value = 0.4614243558677465
$ws.Range("C11").Value = 0.5196793326382828
$ws.Range("D11").Value = 0.5672473917880538
$ws.Range("C12").Value = 22.75735960754645
$ws.Range("D12").Value = 12.29755897919453
$ws.Range("C13").Value = 0.6726814157495205
$ws.Range("D13").Value = 0.7504897251850829
$ws.Range("C14").Value = 1.830998939768056
$ws.Range("D14").Value = 1.667848296399615
$ws.Range("C15").Value = 93.2425410431584
$ws.Range("D15").Value = 22.50803120461897
$ws.Range("C16").Value = -85.23947718154869
$ws.Range("D16").Value = 20.2133753453302
$ws.Range("F16").Value = -101.5175485570292
$ws.Range("G16").Value = -83.14699179957641
$ws.Range("H16").Value = -69.45410721860875
$ws.Range("C17").Value = -77.5746239034152
$ws.Range("D17").Value = 24.91663364700802
$ws.Range("F17").Value = -92.49305820175223
$ws.Range("G17").Value = -72.41392685158225
$ws.Range("H17").Value = -58.53779541063678

Write-Host "Updated statistical description values for rows 2-17 (Mean, STD, and quartile columns)."
